# Update the "Maximum Capacity Factor" workbook:
#  - About!C1 "last updated" date moves forward
#  - MCF sheet: most capacity-factor assumptions are bumped up to 1 (100%)
#  - MCF sheet keeps B17 as the selected / active cell

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date -------------------------
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("C1").Value = (Get-Date -Year 2024 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0)

# --- MCF sheet: capacity factor assumptions -> 1 ------------------------
$mcfWs = $wb.Worksheets.Item("MCF")

$cellsToMax = @("B2", "B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($addr in $cellsToMax) {
    $mcfWs.Range($addr).Value = 1
}

# --- Selection / active sheet state -------------------------------------
$mcfWs.Activate()
$mcfWs.Range("B17").Select()
